{"js": "// Office.js (Word JavaScript API) edit script.\n// Body of: async (context) => { ... }\n//\n// Change 1: Insert a new bulleted list item \"Pass go = 50*properties owned\"\n//           right after the \"Multiple players on same space?\" list item.\n//           The trailing space run that used to end that paragraph moves to\n//           become the trailing space of the new paragraph's text.\n// Change 2: Merge the two runs \"Sell propert\" + \"ies\" into a single run\n//           \"Sell properties\".\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// --- Change 1 ---------------------------------------------------------\n// Find the paragraph that contains \"Multiple players ... space?\"\nlet targetParagraph = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const t = paragraphs.items[i].text;\n  if (t.indexOf(\"Multiple players\") !== -1 && t.indexOf(\"space?\") !== -1) {\n    targetParagraph = paragraphs.items[i];\n    break;\n  }\n}\nif (!targetParagraph) {\n  throw new Error(\"Could not locate the 'Multiple players ... space?' paragraph\");\n}\n\n// That paragraph currently ends with a lone trailing-space run (\" \").\n// Locate and remove it so the paragraph ends right after \"space?\" - the\n// removed space is carried over as the trailing character of the new\n// list item we insert next.\nconst paragraphRange = targetParagraph.getRange();\nconst spaceHits = paragraphRange.search(\" \", { matchCase: true });\nspaceHits.load(\"items\");\nawait context.sync();\n\nif (spaceHits.items.length > 0) {\n  const trailingSpace = spaceHits.items[spaceHits.items.length - 1];\n  trailingSpace.delete();\n  await context.sync();\n}\n\n// Insert the new list item right after the (now trimmed) paragraph. It\n// inherits the \"ListParagraph\" style + numbering (ilvl 0 / numId 1) from\n// its neighbor automatically.\ntargetParagraph.insertParagraph(\"Pass go = 50*properties owned \", \"After\");\nawait context.sync();\n\n// --- Change 2 ---------------------------------------------------------\n// Merge the \"Sell propert\" + \"ies\" runs into a single \"Sell properties\" run.\nconst paragraphs2 = body.paragraphs;\nparagraphs2.load(\"items/text\");\nawait context.sync();\n\nfor (let i = 0; i < paragraphs2.items.length; i++) {\n  if (paragraphs2.items[i].text === \"Sell properties\") {\n    const r = paragraphs2.items[i].getRange();\n    r.insertText(\"Sell properties\", \"Replace\");\n    break;\n  }\n}\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $word.ActiveDocument is the open document.\n#\n# Change 1: Insert a new bulleted list item \"Pass go = 50*properties owned\"\n#           right after the \"Multiple players on same space?\" list item.\n#           The trailing space run that used to end that paragraph moves to\n#           become the trailing space of the new paragraph's text.\n# Change 2: Merge the two runs \"Sell propert\" + \"ies\" into a single run\n#           \"Sell properties\".\n\n$d = $word.ActiveDocument\n\n# --- Change 1 -----------------------------------------------------------\n# Locate \"space?\" (the end of the \"Multiple players on same space?\" item).\n$findRange = $d.Content\n$findRange.Find.Text = \"space?\"\n$found = $findRange.Find.Execute()\nif (-not $found) {\n    throw \"Could not find 'space?' in the document\"\n}\n\n# Split the paragraph right after the matched text - in the original\n# paragraph this leaves only the trailing single-space run behind, now\n# living in a brand-new (but already-ListParagraph/numbered) paragraph.\n$findRange.InsertParagraphAfter()\n\n# Find that brand-new paragraph (its text is just a single space) and type\n# the new bullet text before the leftover space.\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $para = $d.Paragraphs.Item($i)\n    if ($para.Range.Text -eq \" `r\") {\n        $para.Range.InsertBefore(\"Pass go = 50*properties owned\")\n        break\n    }\n}\n\n# --- Change 2 -------------------------------------------------------------\n# Merge the \"Sell propert\" + \"ies\" runs into a single \"Sell properties\" run.\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $para = $d.Paragraphs.Item($i)\n    if ($para.Range.Text -eq \"Sell properties`r\") {\n        $r = $para.Range\n        [void]$r.MoveEnd(1, -1) # exclude the trailing paragraph mark\n        $r.Delete()\n        $r.InsertAfter(\"Sell properties\")\n        break\n    }\n}\n"}
